$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,3).Value = 46063

# Row 3
$ws.Cells.Item(3,3).Value = 46063

# Row 4
$ws.Cells.Item(4,3).Value = 46063

# Row 5
$ws.Cells.Item(5,1).Value = 'A 34341-2024'
$ws.Cells.Item(5,2).Value = 45525
$ws.Cells.Item(5,3).Value = 46063
$ws.Cells.Item(5,6).Value = 'Övriga Aktiebolag'
$ws.Cells.Item(5,7).Value = 14.4
$ws.Cells.Item(5,8).Value = 0
$ws.Cells.Item(5,10).Value = 1
$ws.Cells.Item(5,15).Value = 1
$ws.Cells.Item(5,18).Value = 'Desmeknopp'
$ws.Cells.Item(5,19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 34341-2024 artfynd.xlsx", "A 34341-2024")'
$ws.Cells.Item(5,20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 34341-2024 karta.png", "A 34341-2024")'
$ws.Cells.Item(5,22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 34341-2024 FSC-klagomål.docx", "A 34341-2024")'
$ws.Cells.Item(5,23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 34341-2024 FSC-klagomål mail.docx", "A 34341-2024")'
$ws.Cells.Item(5,24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 34341-2024 tillsynsbegäran.docx", "A 34341-2024")'
$ws.Cells.Item(5,25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 34341-2024 tillsynsbegäran mail.docx", "A 34341-2024")'

# Row 6
$ws.Cells.Item(6,3).Value = 46063

# Row 7
$ws.Cells.Item(7,1).Value = 'A 31213-2023'
$ws.Cells.Item(7,2).Value = 45113
$ws.Cells.Item(7,3).Value = 46063
$ws.Cells.Item(7,7).Value = 6.5
$ws.Cells.Item(7,9).Value = 0
$ws.Cells.Item(7,10).Value = 1
$ws.Cells.Item(7,15).Value = 1
$ws.Cells.Item(7,18).Value = 'Skogsveronika'
$ws.Cells.Item(7,19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 31213-2023 artfynd.xlsx", "A 31213-2023")'
$ws.Cells.Item(7,20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 31213-2023 karta.png", "A 31213-2023")'
$ws.Cells.Item(7,22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 31213-2023 FSC-klagomål.docx", "A 31213-2023")'
$ws.Cells.Item(7,23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 31213-2023 FSC-klagomål mail.docx", "A 31213-2023")'
$ws.Cells.Item(7,24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 31213-2023 tillsynsbegäran.docx", "A 31213-2023")'
$ws.Cells.Item(7,25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 31213-2023 tillsynsbegäran mail.docx", "A 31213-2023")'

# Row 8
$ws.Cells.Item(8,1).Value = 'A 49546-2025'
$ws.Cells.Item(8,2).Value = 45939
$ws.Cells.Item(8,3).Value = 46063
$ws.Cells.Item(8,7).Value = 4.4
$ws.Cells.Item(8,8).Value = 0
$ws.Cells.Item(8,9).Value = 1
$ws.Cells.Item(8,18).Value = 'Igelkottsröksvamp'
$ws.Cells.Item(8,19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 49546-2025 artfynd.xlsx", "A 49546-2025")'
$ws.Cells.Item(8,20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 49546-2025 karta.png", "A 49546-2025")'
$ws.Cells.Item(8,22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 49546-2025 FSC-klagomål.docx", "A 49546-2025")'
$ws.Cells.Item(8,23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 49546-2025 FSC-klagomål mail.docx", "A 49546-2025")'
$ws.Cells.Item(8,24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 49546-2025 tillsynsbegäran.docx", "A 49546-2025")'
$ws.Cells.Item(8,25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 49546-2025 tillsynsbegäran mail.docx", "A 49546-2025")'

# Row 9
$ws.Cells.Item(9,1).Value = 'A 13766-2023'
$ws.Cells.Item(9,2).Value = 45007
$ws.Cells.Item(9,3).Value = 46063
$ws.Cells.Item(9,7).Value = 0.9
$ws.Cells.Item(9,8).Value = 1
$ws.Cells.Item(9,10).Value = 0
$ws.Cells.Item(9,15).Value = 0
$ws.Cells.Item(9,18).Value = 'Större vattensalamander'
$ws.Cells.Item(9,19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 13766-2023 artfynd.xlsx", "A 13766-2023")'
$ws.Cells.Item(9,20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 13766-2023 karta.png", "A 13766-2023")'
$ws.Cells.Item(9,22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 13766-2023 FSC-klagomål.docx", "A 13766-2023")'
$ws.Cells.Item(9,23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 13766-2023 FSC-klagomål mail.docx", "A 13766-2023")'
$ws.Cells.Item(9,24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 13766-2023 tillsynsbegäran.docx", "A 13766-2023")'
$ws.Cells.Item(9,25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 13766-2023 tillsynsbegäran mail.docx", "A 13766-2023")'

# Row 10
$ws.Cells.Item(10,1).Value = 'A 60891-2024'
$ws.Cells.Item(10,2).Value = 45644
$ws.Cells.Item(10,3).Value = 46063
$ws.Cells.Item(10,6).ClearContents()
$ws.Cells.Item(10,7).Value = 16.1
$ws.Cells.Item(10,8).Value = 1
$ws.Cells.Item(10,10).Value = 0
$ws.Cells.Item(10,15).Value = 0
$ws.Cells.Item(10,18).Value = 'Lövgroda'
$ws.Cells.Item(10,19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 60891-2024 artfynd.xlsx", "A 60891-2024")'
$ws.Cells.Item(10,20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 60891-2024 karta.png", "A 60891-2024")'
$ws.Cells.Item(10,22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 60891-2024 FSC-klagomål.docx", "A 60891-2024")'
$ws.Cells.Item(10,23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 60891-2024 FSC-klagomål mail.docx", "A 60891-2024")'
$ws.Cells.Item(10,24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 60891-2024 tillsynsbegäran.docx", "A 60891-2024")'
$ws.Cells.Item(10,25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 60891-2024 tillsynsbegäran mail.docx", "A 60891-2024")'

# Row 11
$ws.Cells.Item(11,3).Value = 46063

# Row 12
$ws.Cells.Item(12,3).Value = 46063

# Row 13
$ws.Cells.Item(13,3).Value = 46063

# Row 14
$ws.Cells.Item(14,1).Value = 'A 40417-2022'
$ws.Cells.Item(14,2).Value = 44823
$ws.Cells.Item(14,3).Value = 46063
$ws.Cells.Item(14,7).Value = 2.3

# Row 15
$ws.Cells.Item(15,1).Value = 'A 34302-2024'
$ws.Cells.Item(15,2).Value = 45524
$ws.Cells.Item(15,3).Value = 46063
$ws.Cells.Item(15,7).Value = 1.9

# Row 16
$ws.Cells.Item(16,1).Value = 'A 60803-2023'
$ws.Cells.Item(16,2).Value = 45260
$ws.Cells.Item(16,3).Value = 46063
$ws.Cells.Item(16,7).Value = 1.6

# Row 17
$ws.Cells.Item(17,1).Value = 'A 32596-2024'
$ws.Cells.Item(17,2).Value = 45513.61667824074
$ws.Cells.Item(17,3).Value = 46063
$ws.Cells.Item(17,7).Value = 2.6

# Row 18
$ws.Cells.Item(18,1).Value = 'A 49536-2025'
$ws.Cells.Item(18,2).Value = 45939.4221875
$ws.Cells.Item(18,3).Value = 46063
$ws.Cells.Item(18,7).Value = 1.5

# Row 19
$ws.Cells.Item(19,1).Value = 'A 49543-2025'
$ws.Cells.Item(19,2).Value = 45939.428622685184
$ws.Cells.Item(19,3).Value = 46063
$ws.Cells.Item(19,7).Value = 1.4

# Row 20
$ws.Cells.Item(20,1).Value = 'A 49549-2025'
$ws.Cells.Item(20,2).Value = 45939
$ws.Cells.Item(20,3).Value = 46063
$ws.Cells.Item(20,7).Value = 0.5

# Row 21
$ws.Cells.Item(21,3).Value = 46063

# Row 22
$ws.Cells.Item(22,1).Value = 'A 18090-2022'
$ws.Cells.Item(22,2).Value = 44684
$ws.Cells.Item(22,3).Value = 46063
$ws.Cells.Item(22,6).ClearContents()
$ws.Cells.Item(22,7).Value = 4.9

# Row 23
$ws.Cells.Item(23,1).Value = 'A 22195-2023'
$ws.Cells.Item(23,2).Value = 45069.74605324074
$ws.Cells.Item(23,3).Value = 46063
$ws.Cells.Item(23,6).Value = 'Övriga Aktiebolag'
$ws.Cells.Item(23,7).Value = 1.1
